# Update the "日志" (log) worksheet with a new day's entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("日志")

# Preserve the current row 12 (C:D) formatting onto the new row 13 first,
# since the new row inherits the wrap-text style that row 12 currently has.
$ws.Range("C12:D12").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122)  # xlPasteFormats

# New row 13: date, time slot and the day's notes (executed == planned).
$ws.Range("A13").Value = 42134
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B13").Value = "17:00-19:00"

$noteText = "1、学习Test类`n2、利用test增加表数据"
$ws.Range("C13").Value = $noteText
$ws.Range("D13").Value = $noteText

# Row 12's C:D cells pick up the other (equivalent) wrap-text style, matching
# how Excel re-indexes styles once a new row with the same style is inserted.
$ws.Range("C2:D2").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)  # xlPasteFormats

# Match the row height Excel computed for the wrapped two-line note.
$ws.Rows.Item(13).RowHeight = 27

# Leave the selection where the user's cursor ended up after entering the row.
$null = $ws.Range("C14").Select()

Write-Host "Updated log sheet with 2015/5/10 entry"
